# Update countries & provincias Spain
# - Refresh the "Datos actualizados..." timestamp (A1)
# - Swap four pairs of countries whose ranking crossed over, and refresh the
#   statistics for every country row whose figures changed in the new data
#   pull (COVID case counters).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Timestamp banner -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 23 de Septiembre de 2020 a las 13:25"

# --- Row 4: Estados Unidos --------------------------------------------
$ws.Range("E4").Value = 2545628
$ws.Range("G4").Value = 20
$ws.Range("H4").Value = 205491

# --- Row 16: Iran -------------------------------------------------------
$ws.Range("B16").Value = 432798
$ws.Range("C16").Value = 3605
$ws.Range("D16").Value = 365846
$ws.Range("E16").Value = 42112
$ws.Range("G16").Value = 184
$ws.Range("H16").Value = 24840

# --- Row 25: Alemania -----------------------------------------------
$ws.Range("B25").Value = 277420
$ws.Range("C25").Value = 244
$ws.Range("E25").Value = 20026
$ws.Range("G25").Value = 3
$ws.Range("H25").Value = 9494

# --- Row 32: Catar --------------------------------------------------
$ws.Range("B32").Value = 124175
$ws.Range("C32").Value = 258
$ws.Range("D32").Value = 121006
$ws.Range("E32").Value = 2957
$ws.Range("G32").Value = 1
$ws.Range("H32").Value = 212

# --- Rows 44/45: Guatemala <-> Emiratos Arabes Unidos swap places, ----
# --- and pick up fresh statistics -------------------------------------
$ws.Range("A44").Value = "Emiratos Arabes Unidos"
$ws.Range("B44").Value = 87530
$ws.Range("C44").Value = 1083
$ws.Range("D44").Value = 76995
$ws.Range("E44").Value = 10129
$ws.Range("G44").Value = 1
$ws.Range("H44").Value = 406

$ws.Range("A45").Value = "Guatemala"
$ws.Range("B45").Value = 86623
$ws.Range("C45").Value = 0
$ws.Range("D45").Value = 75869
$ws.Range("E45").Value = 7617
$ws.Range("G45").Value = 0
$ws.Range("H45").Value = 3137

# --- Rows 54/55: Costa Rica <-> Nepal swap places ---------------------
$ws.Range("A54").Value = "Nepal"
$ws.Range("B54").Value = 67804
$ws.Range("C54").Value = 1172
$ws.Range("D54").Value = 49954
$ws.Range("E54").Value = 17414
$ws.Range("G54").Value = 7
$ws.Range("H54").Value = 436

$ws.Range("A55").Value = "Costa Rica"
$ws.Range("B55").Value = 66689
$ws.Range("C55").Value = 0
$ws.Range("D55").Value = 25706
$ws.Range("E55").Value = 40223
$ws.Range("G55").Value = 0
$ws.Range("H55").Value = 760

# --- Row 65: Ghana ----------------------------------------------------
$ws.Range("B65").Value = 46116
$ws.Range("C65").Value = 54
$ws.Range("D65").Value = 45290
$ws.Range("E65").Value = 529

# --- Row 87: Madagascar ------------------------------------------------
$ws.Range("B87").Value = 16167
$ws.Range("C87").Value = 31
$ws.Range("D87").Value = 14788
$ws.Range("E87").Value = 1153

# --- Row 90: Senegal ----------------------------------------------------
$ws.Range("B90").Value = 14795
$ws.Range("C90").Value = 36
$ws.Range("D90").Value = 11718
$ws.Range("E90").Value = 2774
$ws.Range("G90").Value = 1
$ws.Range("H90").Value = 303

# --- Row 95: Tunez ------------------------------------------------------
$ws.Range("B95").Value = 12479
$ws.Range("C95").Value = 1219
$ws.Range("E95").Value = 9919
$ws.Range("G95").Value = 10
$ws.Range("H95").Value = 174

# --- Row 103: Finlandia --------------------------------------------------
$ws.Range("B103").Value = 9288
$ws.Range("C103").Value = 93
$ws.Range("D103").Value = 7850
$ws.Range("E103").Value = 1095
$ws.Range("G103").Value = 2
$ws.Range("H103").Value = 343

# --- Row 145: Malta -------------------------------------------------------
$ws.Range("B145").Value = 2856
$ws.Range("C145").Value = 42
$ws.Range("D145").Value = 2173
$ws.Range("E145").Value = 658

# --- Rows 148/149: Guyana <-> Islandia swap places -----------------------
$ws.Range("A148").Value = "Islandia"
$ws.Range("B148").Value = 2476
$ws.Range("C148").Value = 57
$ws.Range("D148").Value = 2142
$ws.Range("E148").Value = 324
$ws.Range("H148").Value = 10

$ws.Range("A149").Value = "Guyana"
$ws.Range("B149").Value = 2437
$ws.Range("D149").Value = 1361
$ws.Range("E149").Value = 1009
$ws.Range("H149").Value = 67

# --- Row 168: Vietnam ------------------------------------------------------
$ws.Range("B168").Value = 1069
$ws.Range("C168").Value = 1
$ws.Range("D168").Value = 991
$ws.Range("E168").Value = 43

# --- Rows 214/215: Montserrat <-> Islas Malvinas swap places ---------------
$ws.Range("A214").Value = "Islas Malvinas"
$ws.Range("D214").Value = 13
$ws.Range("H214").Value = 0

$ws.Range("A215").Value = "Montserrat"
$ws.Range("D215").Value = 12
$ws.Range("H215").Value = 1
